$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 554, shifting existing rows 554:668 down to 555:669.
$ws.Rows.Item(554).Insert()

# Populate the newly inserted row 554 with the new data record.
$ws.Cells.Item(554, 1).Value = 3
$ws.Cells.Item(554, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(554, 3).Value = "Coquimbo"
$ws.Cells.Item(554, 4).Value = 45275
$ws.Cells.Item(554, 5).Value = 5
$ws.Cells.Item(554, 6).Value = 100112009
$ws.Cells.Item(554, 7).Value = "Acelga"
$ws.Cells.Item(554, 8).Value = "Sin especificar"
$ws.Cells.Item(554, 9).Value = "Primera"
$ws.Cells.Item(554, 10).Value = 230
$ws.Cells.Item(554, 11).Value = 3000
$ws.Cells.Item(554, 12).Value = 3500
$ws.Cells.Item(554, 13).Value = 3239
$ws.Cells.Item(554, 14).Value = '$/docena de atados (6 kilos)'
$ws.Cells.Item(554, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(554, 16).Value = 540
$ws.Cells.Item(554, 17).Value = 6
$ws.Cells.Item(554, 18).Value = "Hortaliza"
